# LoanStats: drop the "PRODUTO" auto-filter criteria (this also unhides
# every row that the filter had hidden) and restore the plain, unfiltered
# auto-filter range.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoanStats")
$ws.ShowAllData()

# Move the active selection on LoanStats from C18 to C6.
$ws.Range("C6").Select()

# The "verified" picture anchored at row 42 used to stretch down to row 47
# (those rows were hidden by the filter, so it looked tiny); now that the
# rows are visible again it should only span row 42, i.e. 12pt tall. Doing
# this after ShowAllData keeps the anchor's top-left ("from") pinned at
# row 42 like it was before.
$shp = $ws.Shapes.Item(1)
$shp.Height = 12
